$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 05:39"

# --- Reorder "Barbados" to sit right after "San Martin (Parte Holandesa)" (before "Aruba") ---
# This shifts Aruba / Monaco / Seychelles / Islas Turcas y Caicos down by one row (184-188)
# and also refreshes their (and Barbados') daily statistics.
$ws.Range("A184").Value = "Barbados"
$ws.Range("B184").Value = 122
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 98
$ws.Range("E184").Value = 17
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 7

$ws.Range("A185").Value = "Aruba"
$ws.Range("B185").Value = 121
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 105
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 3

$ws.Range("A186").Value = "Monaco"
$ws.Range("B186").Value = 120
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 105
$ws.Range("E186").Value = 11
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 4

$ws.Range("A187").Value = "Seychelles"
$ws.Range("B187").Value = 114
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 39
$ws.Range("E187").Value = 75
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("B188").Value = 114
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 38
$ws.Range("E188").Value = 74
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 2

# --- Reorder "Belice" to sit right after "Polinesia Francesa" (before "San Vicente y las Granadinas") ---
# This shifts San Vicente y las Granadinas / San Martin (Parte Francesa) down by one row (194-196)
# and also refreshes their (and Belice's) daily statistics.
$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 56
$ws.Range("C194").Value = 8
$ws.Range("D194").Value = 30
$ws.Range("E194").Value = 24
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2

$ws.Range("A195").Value = "San Vicente y las Granadinas"
$ws.Range("B195").Value = 54
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 44
$ws.Range("E195").Value = 10
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "San Martin (Parte Francesa)"
$ws.Range("B196").Value = 53
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 41
$ws.Range("E196").Value = 9
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 3

# --- Refresh daily statistics for other countries whose numbers changed ---
# Row 39: Belgica
$ws.Range("B39").Value = 69402
$ws.Range("C39").Value = 651
$ws.Range("D39").Value = 17573
$ws.Range("E39").Value = 41984
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 9845

# Row 51: Honduras
$ws.Range("B51").Value = 42685
$ws.Range("C51").Value = 671
$ws.Range("D51").Value = 5694
$ws.Range("E51").Value = 35623
$ws.Range("G51").Value = 31
$ws.Range("H51").Value = 1368

# Row 172: Mongolia
$ws.Range("B172").Value = 293
$ws.Range("C172").Value = 2
$ws.Range("E172").Value = 63

# Row 176: Camboya
$ws.Range("B176").Value = 240
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 196
$ws.Range("E176").Value = 44

# Row 216: San Bartolome
$ws.Range("B216").Value = 9
$ws.Range("C216").Value = 1
$ws.Range("E216").Value = 3
